# Updated Albert's Computer Folders
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# SourcePath for albert-pc (column D, row 2)
$ws.Range("D2").Value = "C:\Users\Albert\Documents\Princeton\Gregor Lab\Data Analysis\LivemRNA\RawData"

# FISHPath for albert-pc (column D, row 3)
$ws.Range("D3").Value = "C:\Users\Albert\Documents\Princeton\Gregor Lab\Data Analysis\LivemRNA\FISHAnalysisData"

# MS2CodePath for albert-pc (column D, row 7) - newly added
$ws.Range("D7").Value = "C:\Users\Albert\Documents\Princeton\Gregor Lab\Data Analysis\LivemRNA\mRNADynamics"

# Leave the selection where the user ended up after editing column D
[void]$ws.Range("D8").Select()
